$wb = $excel.ActiveWorkbook

# --- Sheet "Overview" (tab 1): File Name / Path And Name / Latest HO Xliff Generate Date ---
$wsOverview = $wb.Worksheets.Item(1)
$wsOverview.Range("A2").Value = "058c0217-89a5-46c8-8cab-28c89641389a.md"
$wsOverview.Range("B2").Value = "e2e\058c0217-89a5-46c8-8cab-28c89641389a.md"
$wsOverview.Range("G2").Value = "2016-09-07 11:21:39"

# Update the existing hyperlink's display text in place (use array conversion,
# not .Item(), so the existing link object is mutated rather than appended).
$overviewLinks = @($wsOverview.Hyperlinks)
$overviewLinks[0].TextToDisplay = "e2e\058c0217-89a5-46c8-8cab-28c89641389a.md"

# --- Sheet "zh-cn" (tab 2): Source File Name / Latest Handoff File / Latest Handoff Datetime ---
$wsZhCn = $wb.Worksheets.Item(2)
$wsZhCn.Range("A2").Value = "058c0217-89a5-46c8-8cab-28c89641389a.md"
$wsZhCn.Range("G2").Value = "058c0217-89a5-46c8-8cab-28c89641389a.dc80d50f3cb31c31d069c5b8a920d85f757c67ff.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-07 11:21:33"

$zhCnLinks = @($wsZhCn.Hyperlinks)
$zhCnLinks[0].TextToDisplay = "058c0217-89a5-46c8-8cab-28c89641389a.md"

# --- Sheet "de-de" (tab 3): Source File Name / Latest Handoff File / Latest Handoff Datetime ---
$wsDeDe = $wb.Worksheets.Item(3)
$wsDeDe.Range("A2").Value = "058c0217-89a5-46c8-8cab-28c89641389a.md"
$wsDeDe.Range("G2").Value = "058c0217-89a5-46c8-8cab-28c89641389a.dc80d50f3cb31c31d069c5b8a920d85f757c67ff.de-de.xlf"
# Same underlying value as the Overview sheet's "Latest HO Xliff Generate Date" (G2)
$wsDeDe.Range("H2").Value = "2016-09-07 11:21:39"

$deDeLinks = @($wsDeDe.Hyperlinks)
$deDeLinks[0].TextToDisplay = "058c0217-89a5-46c8-8cab-28c89641389a.md"
